# Auto update stock data
# Updates the Date_1 column (A) from 2026/01/06 -> 2026/01/07 and refreshes
# the EBITDA column (B) values for each company's latest snapshot row.
#
# Values are written as plain text (matching the workbook's existing
# text-stored numbers/dates) by temporarily forcing a Text number format
# before the assignment, then clearing formatting afterwards so the cell's
# style index is left untouched (same as the source file).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

# row -> (new date, new EBITDA or $null if unchanged)
$updates = @(
    @{ Row = 2;  Date = "2026/01/07"; Ebitda = "7.94" },
    @{ Row = 8;  Date = "2026/01/07"; Ebitda = "8.91" },
    @{ Row = 14; Date = "2026/01/07"; Ebitda = "3.20" },
    @{ Row = 20; Date = "2026/01/07"; Ebitda = "13.13" },
    @{ Row = 26; Date = "2026/01/07"; Ebitda = "11.59" },
    @{ Row = 32; Date = "2026/01/07"; Ebitda = "27.89" },
    @{ Row = 38; Date = "2026/01/07"; Ebitda = $null },
    @{ Row = 44; Date = "2026/01/07"; Ebitda = "13.11" },
    @{ Row = 50; Date = "2026/01/07"; Ebitda = "11.71" },
    @{ Row = 56; Date = "2026/01/07"; Ebitda = "31.99" },
    @{ Row = 62; Date = "2026/01/07"; Ebitda = "11.22" },
    @{ Row = 68; Date = "2026/01/07"; Ebitda = "12.93" },
    @{ Row = 74; Date = "2026/01/07"; Ebitda = "18.22" }
)

foreach ($u in $updates) {
    Set-TextValue ("A" + $u.Row) $u.Date
    if ($null -ne $u.Ebitda) {
        Set-TextValue ("B" + $u.Row) $u.Ebitda
    }
}
